# TestExam.xlsx edit: "Updated test exam again"
#
# - Quote the "Do. Or do not. There is no try." line on the "Task 1" sheet
#   and drop its bold formatting.
# - Change the trailing period to a colon on the "One Ring to rule them
#   all..." line on the "Task 2" sheet and drop its left-indent formatting.
# - Update each sheet's remembered selection and make "Task 2" the active
#   (selected) tab instead of "Task 3".

$wb = $excel.ActiveWorkbook

$wsTask1 = $wb.Worksheets.Item("Task 1")
$wsTask2 = $wb.Worksheets.Item("Task 2")
$wsTask3 = $wb.Worksheets.Item("Task 3")

# --- Text content updates -------------------------------------------------

$wsTask1.Range("C2").Value = '"Do. Or do not. There is no try."'
$wsTask2.Range("C2").Value = '"One Ring to rule them all, One Ring to find them, One Ring to bring them all and in the darkness bind them:"'

# --- Formatting updates -----------------------------------------------------

# "Versions" quotes on Task 1 were bold; clear that back to the default look.
$wsTask1.Range("C2:C6").ClearFormats()

# The first quote on Task 2 had a left-indent alignment; clear it back to
# the default look used by the rest of the column.
$wsTask2.Range("C2").ClearFormats()

# --- View / selection updates ----------------------------------------------
# Order matters: selecting a range on a non-active sheet activates that
# sheet, so "Task 2" is selected last to make it the final active tab.

[void]$wsTask1.Range("C19").Select()
[void]$wsTask3.Range("C14").Select()
[void]$wsTask2.Range("F20").Select()
